# Completed data collection of titres for LB and Kan.
#
# - Row 37 (DG011 #9, LB) had its Volume..Titre columns stored as placeholder text;
#   fill them in for real as numbers.
# - Rows 38-72 are the new Kan-condition titre measurements for DG011, DG012, SLM1042
#   and SLM1043 (9 replicates each).
# - Row 73 (SLM1043 #9, Kan) is the final row and its Volume..Titre values were captured
#   as plain text, the same way the old row 37 was before this edit.
#
# Medians: DG011 #5, DG012 #8, SLM1042 #6 and SLM1042 #8, SLM1043 #1 and SLM1043 #7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 37-72: numeric data (row 37 updated in place, rows 38-72 newly appended).
$numericRows = @(
    @(37, "DG011 #9", "LB", 0.01, 0.000001, 0.00001, 3, 3, 26, 27, 268181818.1818181),
    @(38, "DG011 #1", "Kan", 0.01, 1, 0.1, 37, 40, 5, 2, 3818.181818181818),
    @(39, "DG011 #2", "Kan", 0.01, 1, 0.1, 7, 2, 0, 0, 409.0909090909091),
    @(40, "DG011 #3", "Kan", 0.01, 1, 0.1, 10, 11, 3, 1, 1136.363636363636),
    @(41, "DG011 #4", "Kan", 0.01, 1, 0.1, 12, 8, 2, 2, 1090.909090909091),
    @(42, "DG011 #5", "Kan", 0.01, 1, 0.1, 8, 8, 1, 3, 909.090909090909),
    @(43, "DG011 #6", "Kan", 0.01, 1, 0.1, 2, 2, 1, 2, 318.1818181818181),
    @(44, "DG011 #7", "Kan", 0.01, 1, 0.1, 13, 14, 2, 2, 1409.090909090909),
    @(45, "DG011 #8", "Kan", 0.01, 1, 0.1, 3, 7, 1, 0, 500),
    @(46, "DG011 #9", "Kan", 0.01, 1, 0.1, 4, 8, 1, 0, 590.9090909090909),
    @(47, "DG012 #1", "Kan", 0.01, 1, 0.1, 29, 32, 8, 3, 3272.727272727273),
    @(48, "DG012 #2", "Kan", 0.01, 1, 0.1, 34, 22, 2, 5, 2863.636363636364),
    @(49, "DG012 #3", "Kan", 0.01, 1, 0.1, 23, 18, 1, 2, 2000),
    @(50, "DG012 #4", "Kan", 0.01, 1, 0.1, 26, 24, 9, 9, 3090.909090909091),
    @(51, "DG012 #5", "Kan", 0.01, 1, 0.1, 21, 25, 5, 9, 2727.272727272727),
    @(52, "DG012 #6", "Kan", 0.01, 1, 0.1, 20, 24, 0, 2, 2090.909090909091),
    @(53, "DG012 #7", "Kan", 0.01, 1, 0.1, 17, 13, 2, 3, 1590.909090909091),
    @(54, "DG012 #8", "Kan", 0.01, 1, 0.1, 26, 18, 5, 5, 2454.545454545455),
    @(55, "DG012 #9", "Kan", 0.01, 1, 0.1, 24, 23, 3, 3, 2409.090909090909),
    @(56, "SLM1042 #1", "Kan", 0.01, 1, 0.1, 7, 2, 0, 0, 409.0909090909091),
    @(57, "SLM1042 #2", "Kan", 0.01, 1, 0.1, 8, 11, 1, 2, 1000),
    @(58, "SLM1042 #3", "Kan", 0.01, 1, 0.1, 6, 6, 1, 0, 590.9090909090909),
    @(59, "SLM1042 #4", "Kan", 0.01, 1, 0.1, 7, 4, 1, 1, 590.9090909090909),
    @(60, "SLM1042 #5", "Kan", 0.01, 1, 0.1, 20, 17, 4, 7, 2181.818181818181),
    @(61, "SLM1042 #6", "Kan", 0.01, 1, 0.1, 8, 5, 1, 3, 772.7272727272726),
    @(62, "SLM1042 #7", "Kan", 0.01, 1, 0.1, 9, 11, 5, 2, 1227.272727272727),
    @(63, "SLM1042 #8", "Kan", 0.01, 1, 0.1, 9, 6, 1, 1, 772.7272727272726),
    @(64, "SLM1042 #9", "Kan", 0.01, 1, 0.1, 28, 27, 5, 6, 3000),
    @(65, "SLM1043 #1", "Kan", 0.01, 1, 0.1, 23, 17, 1, 2, 1954.545454545455),
    @(66, "SLM1043 #2", "Kan", 0.01, 1, 0.1, 26, 28, 4, 5, 2863.636363636364),
    @(67, "SLM1043 #3", "Kan", 0.01, 1, 0.1, 15, 18, 3, 3, 1772.727272727272),
    @(68, "SLM1043 #4", "Kan", 0.01, 1, 0.1, 17, 25, 3, 4, 2227.272727272727),
    @(69, "SLM1043 #5", "Kan", 0.01, 1, 0.1, 21, 24, 1, 3, 2227.272727272727),
    @(70, "SLM1043 #6", "Kan", 0.01, 1, 0.1, 15, 21, 1, 4, 1863.636363636364),
    @(71, "SLM1043 #7", "Kan", 0.01, 1, 0.1, 17, 21, 2, 3, 1954.545454545455),
    @(72, "SLM1043 #8", "Kan", 0.01, 1, 0.1, 9, 15, 3, 5, 1454.545454545455),
)

foreach ($dataRow in $numericRows) {
    $r = $dataRow[0]
    $ws.Cells.Item($r, 1).Value = $dataRow[1]
    $ws.Cells.Item($r, 2).Value = $dataRow[2]
    $ws.Cells.Item($r, 3).Value = $dataRow[3]
    $ws.Cells.Item($r, 4).Value = $dataRow[4]
    $ws.Cells.Item($r, 5).Value = $dataRow[5]
    $ws.Cells.Item($r, 6).Value = $dataRow[6]
    $ws.Cells.Item($r, 7).Value = $dataRow[7]
    $ws.Cells.Item($r, 8).Value = $dataRow[8]
    $ws.Cells.Item($r, 9).Value = $dataRow[9]
    $ws.Cells.Item($r, 10).Value = $dataRow[10]
}

# Row 73 (last row): Strain/Condition as normal text.
$ws.Cells.Item(73, 1).Value = "SLM1043 #9"
$ws.Cells.Item(73, 2).Value = "Kan"

# Volume..Titre on row 73 were captured as literal text (not numbers). Force the cells
# to Text format before typing the values so Excel keeps them verbatim (e.g. "10e-1")
# instead of auto-parsing them as numbers, then clear that formatting again so no extra
# styling is left applied to the cells.
$lastRowTextRange = $ws.Range("C73:J73")
$lastRowTextRange.NumberFormat = "@"
$ws.Cells.Item(73, 3).Value = "0.01"
$ws.Cells.Item(73, 4).Value = "10e-1"
$ws.Cells.Item(73, 5).Value = "10e-2"
$ws.Cells.Item(73, 6).Value = "22"
$ws.Cells.Item(73, 7).Value = "24"
$ws.Cells.Item(73, 8).Value = "3"
$ws.Cells.Item(73, 9).Value = "4"
$ws.Cells.Item(73, 10).Value = "2409.090909090909"
$lastRowTextRange.ClearFormats()

